# Updates the "cryptos" price-table sheet: refreshed Price (D) / Volume(1h)
# (E) figures from the latest scrape, plus a couple of rows whose ranking
# swapped places (so Name/Link/Price/Volume all move together).
#
# Price-column values look numeric ("241.59", "0.406", ...) but are stored
# as plain text in the workbook (e.g. "96.401.76" isn't a valid number at
# all - it's a '.'-grouped price). Assigning such a look-alike numeric
# string straight to Range.Value lets Excel's COM layer "helpfully" coerce
# it into a real number (losing formatting like trailing zeros and
# introducing floating-point noise). Forcing the cell to Text format before
# the write - then clearing that format stamp back off afterwards so the
# cell's style index is untouched - keeps the write a plain string, matching
# the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D="96.401.76"; E="  +1.00%  " },
    @{ Row=3; D="3.583.61"; E="  -0.17%  " },
    @{ Row=4; E="  +0.00%  " },
    @{ Row=5; D="241.59"; E="  +1.39%  " },
    @{ Row=6; D="654.82"; E="  +0.14%  " },
    @{ Row=7; D="1.56"; E="  +6.75%  " },
    @{ Row=8; D="0.406"; E="  -0.03%  " },
    @{ Row=9; E="  +0.08%  " },
    @{ Row=10; E="  +3.71%  " },
    @{ Row=11; D="3.580.91"; E="  -0.11%  " },
    @{ Row=12; D="43.24"; E="  +0.62%  " },
    @{ Row=13; E="  +0.75%  " },
    @{ Row=14; D="6.39"; E="  +1.41%  " },
    @{ Row=15; D="4.251.24"; E="  -0.68%  " },
    @{ Row=16; D="96.301.72"; E="  +0.94%  " },
    @{ Row=17; E="  +1.66%  " },
    @{ Row=18; D="3.579.02"; E="  -0.37%  " },
    @{ Row=19; D="7.75"; E="  -5.26%  " },
    @{ Row=20; D="12.53"; E="  -0.16%  " },
    @{ Row=21; D="17.75"; E="  -1.87%  " },
    @{ Row=22; E="  +2.42%  " },
    @{ Row=23; D="511.95"; E="  +0.38%  " },
    @{ Row=24; D="3.43"; E="  -2.49%  " },
    @{ Row=25; D="0.0000203"; E="  +4.03%  " },
    @{ Row=26; E="  +3.21%  " },
    @{ Row=27; D="96.46"; E="  -0.23%  " },
    @{ Row=28; D="12.69"; E="  -0.31%  " },
    @{ Row=29; D="3.777.82"; E="  -0.54%  " },
    @{ Row=30; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.149"; E="  +7.59%  " },
    @{ Row=31; B="PancakeSwap"; C="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D="2.98"; E="  -7.24%  " },
    @{ Row=32; D="11.45"; E="  +1.38%  " },
    @{ Row=33; E="  +0.10%  " },
    @{ Row=34; D="0.183"; E="  +3.56%  " },
    @{ Row=35; D="1.00"; E="  -0.32%  " },
    @{ Row=36; D="31.62"; E="  -0.73%  " },
    @{ Row=37; D="617.35"; E="  +8.34%  " },
    @{ Row=38; B="PolygonEcosystemToken"; C="https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"; D="0.566"; E="  +1.30%  " },
    @{ Row=39; B="RenderToken"; C="https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"; D="8.64"; E="  +4.20%  " },
    @{ Row=40; E="  +9.58%  " },
    @{ Row=41; E="  +0.05%  " },
    @{ Row=42; E="  -0.01%  " },
    @{ Row=43; D="0.907"; E="  -2.21%  " },
    @{ Row=44; E="  +5.78%  " },
    @{ Row=45; D="5.71"; E="  -0.07%  " },
    @{ Row=46; E="  +1.85%  " },
    @{ Row=47; D="34.17"; E="  +1.15%  " },
    @{ Row=48; E="  -0.95%  " },
    @{ Row=49; D="0.0417"; E="  -0.28%  " },
    @{ Row=50; E="  +3.92%  " },
    @{ Row=51; D="3.20"; E="  +2.99%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey('B')) { $ws.Range("B$row").Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Range("C$row").Value = $u.C }

    if ($u.ContainsKey('D')) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.ClearFormats()
    }

    if ($u.ContainsKey('E')) { $ws.Range("E$row").Value = $u.E }
}
